$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# Update Sheet1: B2 value changes from "ezAbezu" to "rreterheh3"
$ws1.Range("B2").Value = "rreterheh3"

# Update Sheet1 selection to L4
$ws1.Activate()
$ws1.Range("L4").Select()

# Sheet2: add new cell I6 with value "ezAbezu", copying style from Sheet1 A2 (style index 2)
$ws2.Range("I6").Value = "ezAbezu"
$ws1.Range("A2").Copy()
$ws2.Range("I6").PasteSpecial(-4122)  # xlPasteFormats

# Sheet2 selection to I6
$ws2.Activate()
$ws2.Range("I6").Select()

$ws1.Activate()
